$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between rows 2-3 and rows 4-5 (columns D, M, N, O, P, S).
# Row 2 <-> Row 4, Row 3 <-> Row 5

$ws.Range("D2").Value = 44574
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6500
$ws.Range("S2").Value = 3250

$ws.Range("D3").Value = 44574
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 5000
$ws.Range("O3").Value = 5000
$ws.Range("P3").Value = 5000
$ws.Range("S3").Value = 2500

$ws.Range("D4").Value = 44223
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 3500
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3750
$ws.Range("S4").Value = 1875

$ws.Range("D5").Value = 44223
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 3000
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("S5").Value = 1500
